# Apply the Zalera_Profits.xlsx scheduled-runner update.
# The sheet tracks FFXIV Leve crafting economics; columns H:N hold the
# current market-board snapshot (currentAveragePrice / currentAveragePriceNQ/HQ)
# and the resulting Leve price/profit figures (LevePriceNQ/HQ, LeveProfitNQ/HQ).
# This refresh overwrites those market-derived columns per row with the latest
# values across all eight job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(51, 8).Value = 54985.145
$ws.Cells.Item(51, 9).Value = 7073.25
$ws.Cells.Item(51, 10).Value = 74149.89999999999
$ws.Cells.Item(51, 11).Value = 7073.25
$ws.Cells.Item(51, 12).Value = 74149.89999999999
$ws.Cells.Item(51, 13).Value = -6589.25
$ws.Cells.Item(51, 14).Value = -75117.89999999999
$ws.Cells.Item(74, 8).Value = 6995.8
$ws.Cells.Item(74, 9).Value = 5493.7
$ws.Cells.Item(74, 11).Value = 5493.7
$ws.Cells.Item(74, 13).Value = -4557.7
$ws.Cells.Item(77, 8).Value = 6995.8
$ws.Cells.Item(77, 9).Value = 5493.7
$ws.Cells.Item(77, 11).Value = 27468.5
$ws.Cells.Item(77, 13).Value = -22788.5
$ws.Cells.Item(96, 8).Value = 6255583.5
$ws.Cells.Item(96, 9).Value = 249.83333
$ws.Cells.Item(96, 11).Value = 749.49999
$ws.Cells.Item(96, 13).Value = 623.50001
$ws.Cells.Item(112, 8).Value = 2653.0557
$ws.Cells.Item(112, 9).Value = 1860
$ws.Cells.Item(112, 11).Value = 5580
$ws.Cells.Item(112, 13).Value = -4472
$ws.Cells.Item(132, 8).Value = 1487.6
$ws.Cells.Item(132, 9).Value = 1254.9231
$ws.Cells.Item(132, 11).Value = 3764.7693
$ws.Cells.Item(132, 13).Value = -1234.7693
$ws.Cells.Item(137, 8).Value = 41683150
$ws.Cells.Item(137, 10).Value = 19782
$ws.Cells.Item(137, 12).Value = 59346
$ws.Cells.Item(137, 14).Value = -64446
$ws.Cells.Item(141, 8).Value = 2638
$ws.Cells.Item(141, 9).Value = 1850.6666
$ws.Cells.Item(141, 11).Value = 5551.9998
$ws.Cells.Item(141, 13).Value = -371.9997999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(25, 8).Value = 515.5
$ws.Cells.Item(25, 9).Value = 458.6
$ws.Cells.Item(25, 11).Value = 458.6
$ws.Cells.Item(25, 13).Value = -56.60000000000002
$ws.Cells.Item(45, 8).Value = 1316.1428
$ws.Cells.Item(45, 9).Value = 1340.4615
$ws.Cells.Item(45, 11).Value = 1340.4615
$ws.Cells.Item(45, 13).Value = -963.4614999999999
$ws.Cells.Item(61, 8).Value = 4595.4287
$ws.Cells.Item(61, 9).Value = 3465.5417
$ws.Cells.Item(61, 11).Value = 3465.5417
$ws.Cells.Item(61, 13).Value = -3253.5417
$ws.Cells.Item(74, 8).Value = 4350.8335
$ws.Cells.Item(74, 9).Value = 2231.7144
$ws.Cells.Item(74, 11).Value = 2231.7144
$ws.Cells.Item(74, 13).Value = -1357.7144
$ws.Cells.Item(77, 8).Value = 4350.8335
$ws.Cells.Item(77, 9).Value = 2231.7144
$ws.Cells.Item(77, 11).Value = 11158.572
$ws.Cells.Item(77, 13).Value = -6790.572
$ws.Cells.Item(82, 8).Value = 36393.668
$ws.Cells.Item(82, 9).Value = 30000
$ws.Cells.Item(82, 10).Value = 39590.5
$ws.Cells.Item(82, 11).Value = 30000
$ws.Cells.Item(82, 12).Value = 39590.5
$ws.Cells.Item(82, 13).Value = -29639
$ws.Cells.Item(82, 14).Value = -40312.5
$ws.Cells.Item(85, 8).Value = 36393.668
$ws.Cells.Item(85, 9).Value = 30000
$ws.Cells.Item(85, 10).Value = 39590.5
$ws.Cells.Item(85, 11).Value = 30000
$ws.Cells.Item(85, 12).Value = 39590.5
$ws.Cells.Item(85, 13).Value = -28752
$ws.Cells.Item(85, 14).Value = -42086.5
$ws.Cells.Item(110, 8).Value = 85556800
$ws.Cells.Item(110, 9).Value = 96251336
$ws.Cells.Item(110, 10).Value = 500
$ws.Cells.Item(110, 11).Value = 96251336
$ws.Cells.Item(110, 12).Value = 500
$ws.Cells.Item(110, 13).Value = -96249291
$ws.Cells.Item(110, 14).Value = -4590
$ws.Cells.Item(132, 8).Value = 2330.309
$ws.Cells.Item(132, 9).Value = 1743.7609
$ws.Cells.Item(132, 11).Value = 5231.2827
$ws.Cells.Item(132, 13).Value = -2701.2827
$ws.Cells.Item(136, 8).Value = 4595.4287
$ws.Cells.Item(136, 9).Value = 3465.5417
$ws.Cells.Item(136, 11).Value = 10396.6251
$ws.Cells.Item(136, 13).Value = -7846.625100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 3009
$ws.Cells.Item(94, 9).Value = 3009
$ws.Cells.Item(94, 11).Value = 3009
$ws.Cells.Item(94, 13).Value = -2558
$ws.Cells.Item(99, 8).Value = 2917
$ws.Cells.Item(99, 9).Value = 2563.5334
$ws.Cells.Item(99, 10).Value = 3977.4
$ws.Cells.Item(99, 11).Value = 2563.5334
$ws.Cells.Item(99, 12).Value = 3977.4
$ws.Cells.Item(99, 13).Value = -1065.5334
$ws.Cells.Item(99, 14).Value = -6973.4
$ws.Cells.Item(102, 8).Value = 25580.666
$ws.Cells.Item(102, 9).Value = 11956.8
$ws.Cells.Item(102, 11).Value = 11956.8
$ws.Cells.Item(102, 13).Value = -8711.799999999999
$ws.Cells.Item(134, 8).Value = 5187.2354
$ws.Cells.Item(134, 9).Value = 2726.182
$ws.Cells.Item(134, 10).Value = 9699.166999999999
$ws.Cells.Item(134, 11).Value = 8178.545999999999
$ws.Cells.Item(134, 12).Value = 29097.501
$ws.Cells.Item(134, 13).Value = -5643.545999999999
$ws.Cells.Item(134, 14).Value = -34167.501

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 732.93335
$ws.Cells.Item(22, 9).Value = 702.9091
$ws.Cells.Item(22, 10).Value = 815.5
$ws.Cells.Item(22, 11).Value = 702.9091
$ws.Cells.Item(22, 12).Value = 815.5
$ws.Cells.Item(22, 13).Value = -352.9091
$ws.Cells.Item(22, 14).Value = -1515.5
$ws.Cells.Item(31, 8).Value = 40003720
$ws.Cells.Item(31, 9).Value = 100000920
$ws.Cells.Item(31, 11).Value = 100000920
$ws.Cells.Item(31, 13).Value = -100000625
$ws.Cells.Item(34, 8).Value = 40003720
$ws.Cells.Item(34, 9).Value = 100000920
$ws.Cells.Item(34, 11).Value = 100000920
$ws.Cells.Item(34, 13).Value = -100000718
$ws.Cells.Item(58, 8).Value = 4848.5884
$ws.Cells.Item(58, 9).Value = 3274
$ws.Cells.Item(58, 10).Value = 6248.222
$ws.Cells.Item(58, 11).Value = 3274
$ws.Cells.Item(58, 12).Value = 6248.222
$ws.Cells.Item(58, 13).Value = -3071
$ws.Cells.Item(58, 14).Value = -6654.222
$ws.Cells.Item(80, 8).Value = 57482.25
$ws.Cells.Item(80, 10).Value = 57482.25
$ws.Cells.Item(80, 12).Value = 57482.25
$ws.Cells.Item(80, 14).Value = -59728.25
$ws.Cells.Item(83, 8).Value = 57482.25
$ws.Cells.Item(83, 10).Value = 57482.25
$ws.Cells.Item(83, 12).Value = 172446.75
$ws.Cells.Item(83, 14).Value = -183678.75
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 14).ClearContents()
$ws.Cells.Item(99, 8).Value = 5413.5
$ws.Cells.Item(99, 9).Value = 4123.5
$ws.Cells.Item(99, 10).Value = 7993.5
$ws.Cells.Item(99, 11).Value = 4123.5
$ws.Cells.Item(99, 12).Value = 7993.5
$ws.Cells.Item(99, 13).Value = -2625.5
$ws.Cells.Item(99, 14).Value = -10989.5
$ws.Cells.Item(126, 8).Value = 5413.5
$ws.Cells.Item(126, 9).Value = 4123.5
$ws.Cells.Item(126, 10).Value = 7993.5
$ws.Cells.Item(126, 11).Value = 12370.5
$ws.Cells.Item(126, 12).Value = 23980.5
$ws.Cells.Item(126, 13).Value = -9900.5
$ws.Cells.Item(126, 14).Value = -28920.5
$ws.Cells.Item(134, 8).Value = 6510.5
$ws.Cells.Item(134, 9).Value = 5805.6313
$ws.Cells.Item(134, 10).Value = 9189
$ws.Cells.Item(134, 11).Value = 17416.8939
$ws.Cells.Item(134, 12).Value = 27567
$ws.Cells.Item(134, 13).Value = -14881.8939
$ws.Cells.Item(134, 14).Value = -32637
$ws.Cells.Item(136, 8).Value = 4848.5884
$ws.Cells.Item(136, 9).Value = 3274
$ws.Cells.Item(136, 10).Value = 6248.222
$ws.Cells.Item(136, 11).Value = 9822
$ws.Cells.Item(136, 12).Value = 18744.666
$ws.Cells.Item(136, 13).Value = -7272
$ws.Cells.Item(136, 14).Value = -23844.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 2200.8125
$ws.Cells.Item(2, 9).Value = 19.5
$ws.Cells.Item(2, 10).Value = 3897.389
$ws.Cells.Item(2, 11).Value = 117
$ws.Cells.Item(2, 12).Value = 23384.334
$ws.Cells.Item(2, 13).Value = -4
$ws.Cells.Item(2, 14).Value = -23610.334
$ws.Cells.Item(3, 8).Value = 5712.5
$ws.Cells.Item(3, 9).Value = 5712.5
$ws.Cells.Item(3, 11).Value = 17137.5
$ws.Cells.Item(3, 13).Value = -17025.5
$ws.Cells.Item(81, 8).Value = 3791.4443
$ws.Cells.Item(81, 10).Value = 10015
$ws.Cells.Item(81, 12).Value = 30045
$ws.Cells.Item(81, 14).Value = -32291
$ws.Cells.Item(84, 8).Value = 3791.4443
$ws.Cells.Item(84, 10).Value = 10015
$ws.Cells.Item(84, 12).Value = 90135
$ws.Cells.Item(84, 14).Value = -101367
$ws.Cells.Item(136, 8).Value = 2486.6667
$ws.Cells.Item(136, 9).Value = 2486.6667
$ws.Cells.Item(136, 11).Value = 7460.000100000001
$ws.Cells.Item(136, 13).Value = -2360.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(99, 8).Value = 13268.2
$ws.Cells.Item(99, 9).Value = 13268.2
$ws.Cells.Item(99, 11).Value = 13268.2
$ws.Cells.Item(99, 13).Value = -11022.2
$ws.Cells.Item(122, 8).Value = 5005.25
$ws.Cells.Item(122, 9).Value = 4419.9
$ws.Cells.Item(122, 10).Value = 6468.625
$ws.Cells.Item(122, 11).Value = 13259.7
$ws.Cells.Item(122, 12).Value = 19405.875
$ws.Cells.Item(122, 13).Value = -10809.7
$ws.Cells.Item(122, 14).Value = -24305.875
$ws.Cells.Item(132, 8).Value = 5748.4736
$ws.Cells.Item(132, 9).Value = 3029.3
$ws.Cells.Item(132, 10).Value = 8769.777
$ws.Cells.Item(132, 11).Value = 9087.900000000001
$ws.Cells.Item(132, 12).Value = 26309.331
$ws.Cells.Item(132, 13).Value = -6557.900000000001
$ws.Cells.Item(132, 14).Value = -31369.331

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 11673.909
$ws.Cells.Item(40, 9).Value = 12191.3
$ws.Cells.Item(40, 11).Value = 12191.3
$ws.Cells.Item(40, 13).Value = -12055.3
$ws.Cells.Item(46, 8).Value = 4384.077
$ws.Cells.Item(46, 9).Value = 1892.3077
$ws.Cells.Item(46, 11).Value = 1892.3077
$ws.Cells.Item(46, 13).Value = -1704.3077
$ws.Cells.Item(68, 8).Value = 2794.1
$ws.Cells.Item(68, 9).Value = 2588.8
$ws.Cells.Item(68, 11).Value = 2588.8
$ws.Cells.Item(68, 13).Value = -1839.8
$ws.Cells.Item(71, 8).Value = 2794.1
$ws.Cells.Item(71, 9).Value = 2588.8
$ws.Cells.Item(71, 11).Value = 12944
$ws.Cells.Item(71, 13).Value = -9200
